# This workbook row-block got re-ordered upstream: the data that used to
# live on one row now belongs on another (same sheet, same columns), while
# the row numbers themselves stay put. We therefore read each row's values
# first (into variables) and only then write them back in the new order,
# so that reads are never affected by earlier writes.
#
# Columns Y and AA ("Startdatum"/"Slutdatum") hold the literal text
# "2026-02-07" in every one of these rows, so they are identical before
# and after the rotation; we deliberately leave them untouched (A:X and
# AB:AY only) to avoid Excel's automatic text->date coercion when values
# are written back through Value2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- snapshot every row involved, before any writes happen ---
$r8a  = $ws.Range("A8:X8").Value2
$r8b  = $ws.Range("AB8:AY8").Value2
$r9a  = $ws.Range("A9:X9").Value2
$r9b  = $ws.Range("AB9:AY9").Value2

$r10a = $ws.Range("A10:X10").Value2
$r10b = $ws.Range("AB10:AY10").Value2
$r11a = $ws.Range("A11:X11").Value2
$r11b = $ws.Range("AB11:AY11").Value2
$r12a = $ws.Range("A12:X12").Value2
$r12b = $ws.Range("AB12:AY12").Value2
$r13a = $ws.Range("A13:X13").Value2
$r13b = $ws.Range("AB13:AY13").Value2

$r31a = $ws.Range("A31:X31").Value2
$r31b = $ws.Range("AB31:AY31").Value2
$r32a = $ws.Range("A32:X32").Value2
$r32b = $ws.Range("AB32:AY32").Value2
$r33a = $ws.Range("A33:X33").Value2
$r33b = $ws.Range("AB33:AY33").Value2
$r34a = $ws.Range("A34:X34").Value2
$r34b = $ws.Range("AB34:AY34").Value2

# --- rows 8 and 9 swap their full content ---
$ws.Range("A8:X8").Value2   = $r9a
$ws.Range("AB8:AY8").Value2 = $r9b
$ws.Range("A9:X9").Value2   = $r8a
$ws.Range("AB9:AY9").Value2 = $r8b

# --- rows 10-13 rotate: 10<-11, 11<-12, 12<-13, 13<-10 ---
$ws.Range("A10:X10").Value2   = $r11a
$ws.Range("AB10:AY10").Value2 = $r11b
$ws.Range("A11:X11").Value2   = $r12a
$ws.Range("AB11:AY11").Value2 = $r12b
$ws.Range("A12:X12").Value2   = $r13a
$ws.Range("AB12:AY12").Value2 = $r13b
$ws.Range("A13:X13").Value2   = $r10a
$ws.Range("AB13:AY13").Value2 = $r10b

# --- rows 31-34 rotate: 31<-34, 32<-31, 33<-32, 34<-33 ---
$ws.Range("A31:X31").Value2   = $r34a
$ws.Range("AB31:AY31").Value2 = $r34b
$ws.Range("A32:X32").Value2   = $r31a
$ws.Range("AB32:AY32").Value2 = $r31b
$ws.Range("A33:X33").Value2   = $r32a
$ws.Range("AB33:AY33").Value2 = $r32b
$ws.Range("A34:X34").Value2   = $r33a
$ws.Range("AB34:AY34").Value2 = $r33b
